# Generate Report for Handoff
# Updates the "b.md" row on the Overview, zh-cn and de-de sheets to reflect
# that the handoff xliff files were (re)generated, per the latest report run.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b223eb013c257fb9e2351bf965bcf4f75fa51b97/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4393e1738923a7c1269a7a14cad6abf53a102a7e/e2e/b.md."

# ---------------------------------------------------------------------
# Overview sheet - row 3 is the b.md entry
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-22 04:47:56"

# ---------------------------------------------------------------------
# zh-cn sheet - row 3 is the b.md entry
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("F3").Value = "False"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-22 04:47:52"
$zhcn.Range("P3").Value = $errorDetail
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet - row 3 is the b.md entry
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "False"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-22 04:47:56"
$dede.Range("P3").Value = $errorDetail
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
